# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet, pushing the existing "Late" / heading / "Outstanding" columns one
# place to the right, then make that sheet the active sheet/tab with the
# selection parked on T6 (matching the saved view state in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column at N (shifts N:P -> O:Q).
$ws.Columns("N").Insert()

# Give the freshly inserted column the same width as its neighbours used
# to have before the shift (column M / former N were ~10.71 characters
# wide).
$ws.Columns("N").ColumnWidth = 10.7109375

# Make "Repayment schedule" the active sheet/tab and park the selection on
# T6, as captured in the workbook's saved view state.
$ws.Activate()
$ws.Range("T6").Select()
